# Updated cryptos list on Tue May 23 16:14:29 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.357.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "'1.860.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").Value = "'315.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").Value = "'0.4616"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D8").Value = "'0.3718"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").Value = "'0.07307"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("D10").Value = "'0.8914"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.09%  "

$ws.Range("D11").Value = "'20.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").Value = "'0.07860"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.399"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.796.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.97%  "

$ws.Range("D15").Value = "'6.549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").Value = "'91.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "'1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "'0.000008933"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "

$ws.Range("D20").Value = "'14.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").Value = "'27.387.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "'5.137"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").Value = "'10.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'2.107.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("D25").Value = "'1.931"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.35%  "

$ws.Range("D26").Value = "'152.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").Value = "'18.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.094"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'116.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("D31").Value = "'0.08835"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7736"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.46%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.063"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("D34").Value = "'1.174"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.38%  "

$ws.Range("D35").Value = "'4.527"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.79%  "

$ws.Range("D36").Value = "'2.696"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.05%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.01963"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.43%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.079"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").Value = "'0.05268"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("D40").Value = "'2.970"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("D41").Value = "'7.065"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("D42").Value = "'0.5144"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").Value = "'0.1645"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("D44").Value = "'8.420"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("D45").Value = "'0.4813"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'10.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.57%  "

$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").Value = "'103.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").Value = "'1.647"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.25%  "

$ws.Range("D50").Value = "'0.06225"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("D51").Value = "'66.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "

